$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) PROFESSIONAL SUMMARY paragraph: "affecting all Black and Asian-American
#    voters, developed" -> "affecting 50M voters, developed"
# ---------------------------------------------------------------------------
$summary = $d.Paragraphs.Item(4).Range
$f = $summary.Find
$f.ClearFormatting()
$f.Execute("affecting all Black and Asian-American voters, developed", $true, $false, $false, $false, $false, `
           $true, 1, $false, "affecting 50M voters, developed", 2)

# ---------------------------------------------------------------------------
# 2) KEY ACHIEVEMENTS AND IMPACT section: replace the 4 bullet paragraphs
#    with 6 new bullet paragraphs.
# ---------------------------------------------------------------------------
# Locate the "Impact" Heading3 paragraph that introduces the bullet list
# (it is the paragraph right after "KEY ACHIEVEMENTS AND IMPACT").
$achHeading = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($txt -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $achHeading = $i
        break
    }
}
$impactIdx = $achHeading + 1
$firstBulletIdx = $impactIdx + 1

# Anchor the six new (plain-style) paragraphs on the first existing bullet
# paragraph so they do not inherit any heading style, then fill them in,
# then delete the four old bullet paragraphs.
$anchor = $d.Paragraphs.Item($firstBulletIdx).Range
$anchor.InsertParagraphAfter()
$d.Paragraphs.Item($firstBulletIdx + 1).Range.InsertParagraphAfter()
$d.Paragraphs.Item($firstBulletIdx + 2).Range.InsertParagraphAfter()
$d.Paragraphs.Item($firstBulletIdx + 3).Range.InsertParagraphAfter()
$d.Paragraphs.Item($firstBulletIdx + 4).Range.InsertParagraphAfter()
$d.Paragraphs.Item($firstBulletIdx + 5).Range.InsertParagraphAfter()

$n1 = $firstBulletIdx + 1
$n2 = $firstBulletIdx + 2
$n3 = $firstBulletIdx + 3
$n4 = $firstBulletIdx + 4
$n5 = $firstBulletIdx + 5
$n6 = $firstBulletIdx + 6

$bullet = [char]0x2022

# --- New paragraph 1: algorithmic innovation / 73.5% ---
$d.Paragraphs.Item($n1).Range.Text = "$bullet Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs **73.5%**"
$r = $d.Paragraphs.Item($n1).Range
$fr = $r.Find
$fr.ClearFormatting()
$fr.Execute("73.5%", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Font.Bold = 1
$r.Font.Color = 5258796

# --- New paragraph 2: $4.7M savings enabled nonprofit access ---
$d.Paragraphs.Item($n2).Range.Text = "$bullet **`$4.7M** savings enabled nonprofit access"
$r = $d.Paragraphs.Item($n2).Range
$fr = $r.Find
$fr.ClearFormatting()
$fr.Execute("`$4.7M", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Font.Bold = 1
$r.Font.Color = 5258796

# --- New paragraph 3: Legal precedent ---
$d.Paragraphs.Item($n3).Range.Text = "$bullet Legal precedent: Data analysis utilized in Supreme Court case"

# --- New paragraph 4: Expert methodology ---
$d.Paragraphs.Item($n4).Range.Text = "$bullet Expert methodology validated at highest judicial level"

# --- New paragraph 5: Breakthrough demographic discovery ---
$d.Paragraphs.Item($n5).Range.Text = "$bullet Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"

# --- New paragraph 6: 178% accuracy improvement ---
$d.Paragraphs.Item($n6).Range.Text = "$bullet **178%** accuracy improvement in racial classification algorithms"
$r = $d.Paragraphs.Item($n6).Range
$fr = $r.Find
$fr.ClearFormatting()
$fr.Execute("178%", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Font.Bold = 1
$r.Font.Color = 5258796

# Delete the 4 original bullet paragraphs: the first one sits right before
# our 6 new paragraphs, the remaining three sit right after them.
$d.Paragraphs.Item($firstBulletIdx).Range.Delete()
$d.Paragraphs.Item($n6).Range.Delete()
$d.Paragraphs.Item($n6).Range.Delete()
$d.Paragraphs.Item($n6).Range.Delete()

Write-Output "Done"
